# Applies the recorded edits: updated student/record ID numbers in column C
# (rows 2-10) and a refreshed cell selection on Sheet1, matching the
# "File upload feature modification + added more data" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Updated ID values in column C (rows 2-10) ---
$ws.Range("C2").Value = 201003
$ws.Range("C3").Value = 201004
$ws.Range("C4").Value = 201005
$ws.Range("C5").Value = 201006
$ws.Range("C6").Value = 201007
$ws.Range("C7").Value = 201008
$ws.Range("C8").Value = 201009
$ws.Range("C9").Value = 201010
$ws.Range("C10").Value = 201011

# --- Refresh the active selection / view on the sheet ---
$ws.Activate() | Out-Null
$ws.Range("J20").Select() | Out-Null
